# Apply hierarchical numbering prefixes to section header cells in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "1. Test automation requirements for IFXX"
$ws.Range("A2").Value = "1.1. Waaalo"
$ws.Range("A3").Value = "1.2. Pre-procession of test automation"
$ws.Range("A4").Value = "1.2.1. Test procedure edition"
$ws.Range("A9").Value = "1.2.2. Test properties"
$ws.Range("A13").Value = "1.2.3. Automated interactions"
$ws.Range("A18").Value = "1.3. procession of test automation"
$ws.Range("A19").Value = "1.3.1. oho ya oho"
$ws.Range("A26").Value = "1.3.2. zbob"
$ws.Range("A27").Value = "1.3.3. Oho yaa"
